$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("D2").Value = "62.937.32"
$ws.Range("E2").Value = "  +2.10%  "
$ws.Range("D3").Value = "3.034.57"
$ws.Range("E3").Value = "  +1.27%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").Value = "'593.72"
$ws.Range("E5").Value = "  -0.59%  "
$ws.Range("D6").Value = "'154.19"
$ws.Range("E6").Value = "  +7.17%  "
$ws.Range("E7").Value = "  -0.05%  "
$ws.Range("D8").Value = "3.030.20"
$ws.Range("E8").Value = "  +1.15%  "
$ws.Range("D9").Value = "'0.518"
$ws.Range("E9").Value = "  -0.78%  "
$ws.Range("E10").Value = "  +8.50%  "
$ws.Range("D11").Value = "'0.151"
$ws.Range("E11").Value = "  +1.96%  "
$ws.Range("E12").Value = "  +0.67%  "
$ws.Range("E13").Value = "  +2.28%  "
$ws.Range("D14").Value = "'35.61"
$ws.Range("E14").Value = "  +3.93%  "
$ws.Range("E15").Value = "  +2.25%  "
$ws.Range("D16").Value = "3.530.60"
$ws.Range("E16").Value = "  +1.06%  "
$ws.Range("D17").Value = "'7.11"
$ws.Range("E17").Value = "  +1.05%  "
$ws.Range("D18").Value = "62.877.84"
$ws.Range("D19").Value = "3.035.08"
$ws.Range("E19").Value = "  +1.14%  "
$ws.Range("D20").Value = "'452.40"
$ws.Range("E20").Value = "  -0.44%  "
$ws.Range("D21").Value = "'14.35"
$ws.Range("E21").Value = "  +2.52%  "
$ws.Range("D22").Value = "'0.696"
$ws.Range("E22").Value = "  +1.28%  "
$ws.Range("D23").Value = "'7.48"
$ws.Range("E23").Value = "  +1.59%  "
$ws.Range("D24").Value = "'83.07"
$ws.Range("E24").Value = "  +1.12%  "
$ws.Range("D25").Value = "'2.31"
$ws.Range("E25").Value = "  +3.44%  "
$ws.Range("D26").Value = "'11.11"
$ws.Range("E26").Value = "  +5.59%  "
$ws.Range("D27").Value = "'12.35"
$ws.Range("E27").Value = "  +1.64%  "
$ws.Range("D29").Value = "'7.48"
$ws.Range("E29").Value = "  +5.83%  "
$ws.Range("E30").Value = "  +0.65%  "
$ws.Range("E31").Value = "  -0.06%  "
$ws.Range("D32").Value = "'2.21"
$ws.Range("E32").Value = "  +6.20%  "
$ws.Range("D33").Value = "'27.60"
$ws.Range("E33").Value = "  +0.19%  "
$ws.Range("E34").Value = "  +1.79%  "
$ws.Range("D35").Value = "0.0₃0871"
$ws.Range("E35").Value = "  +5.50%  "
$ws.Range("E36").Value = "  +1.88%  "
$ws.Range("D37").Value = "'5.93"
$ws.Range("E37").Value = "  +2.95%  "
$ws.Range("D38").Value = "'3.20"
$ws.Range("E38").Value = "  +10.66%  "
$ws.Range("E39").Value = "  +1.96%  "
$ws.Range("B40").Value = "Kaspa"
$ws.Range("C40").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D40").Value = "'0.129"
$ws.Range("E40").Value = "  +5.84%  "
$ws.Range("B41").Value = "OKB"
$ws.Range("C41").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D41").Value = "'50.62"
$ws.Range("E41").Value = "  +0.54%  "
$ws.Range("E42").Value = "  -1.91%  "
$ws.Range("D43").Value = "'0.307"
$ws.Range("E43").Value = "  +14.56%  "
$ws.Range("D44").Value = "'41.47"
$ws.Range("E44").Value = "  +6.95%  "
$ws.Range("D45").Value = "'394.79"
$ws.Range("E45").Value = "  -0.26%  "
$ws.Range("E46").Value = "  +1.48%  "
$ws.Range("D47").Value = "2.731.51"
$ws.Range("E47").Value = "  +0.37%  "
$ws.Range("D48").Value = "'132.70"
$ws.Range("E48").Value = "  -0.77%  "
$ws.Range("D50").Value = "'2.24"
$ws.Range("D51").Value = "'24.38"
$ws.Range("E51").Value = "  +3.74%  "
